# Updating the Staging testdata
#
# Column K ("ExpectedFilenames") on Sheet1 held a mixed Pfizer/Takeda list of
# expected report filenames. Replace it with the refreshed Takeda-only list
# (new "Standard..." naming convention) and drop the now-unused tail rows
# (K14:K18) that used to hold the trailing Pfizer entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value  = "StandardExcelReport-Takeda - MM Maintenance-Clinical-2023_"
$ws.Range("K3").Value  = "ExcelReport-Takeda-MM Maintenance-Clinical-"
$ws.Range("K4").Value  = "WordReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("K5").Value  = "StandardExcelReport-Takeda - MM Maintenance-Economic-2023_"
$ws.Range("K6").Value  = "ExcelReport-Takeda-MM Maintenance-Economic-"
$ws.Range("K7").Value  = "WordReport-Takeda - MM Maintenance-Economic-"
$ws.Range("K8").Value  = "StandardExcelReport-Takeda - MM Maintenance-Quality of Life-2023_"
$ws.Range("K9").Value  = "ExcelReport-Takeda-MM Maintenance-Quality of Life-"
$ws.Range("K10").Value = "WordReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("K11").Value = "StandardExcelReport-Takeda - MM Maintenance-Real-world Evidence-2023_"
$ws.Range("K12").Value = "ExcelReport-Takeda-MM Maintenance-Real-world Evidence-"
$ws.Range("K13").Value = "WordReport-Takeda - MM Maintenance-Real-world Evidence-"

# The old list ran through K18; the refreshed list stops at K13, so clear
# the now-empty tail instead of leaving stale Pfizer values behind.
$ws.Range("K14:K18").ClearContents()

# Reflect the updated viewport/selection recorded for the sheet.
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("K11").Select() | Out-Null
